# Auto-generated edit script: update market price / profit figures
# per Sheets/Adamantoise_Profits.xlsx diff (scheduled runner sync)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6374.24
$ws.Range("I62").Value = 6026.353
$ws.Range("J62").Value = 7113.5
$ws.Range("K62").Value = 6026.353
$ws.Range("L62").Value = 7113.5
$ws.Range("M62").Value = -5402.353
$ws.Range("N62").Value = -8361.5
$ws.Range("H65").Value = 6374.24
$ws.Range("I65").Value = 6026.353
$ws.Range("J65").Value = 7113.5
$ws.Range("K65").Value = 30131.765
$ws.Range("L65").Value = 35567.5
$ws.Range("M65").Value = -27011.765
$ws.Range("N65").Value = -41807.5
$ws.Range("H98").Value = 1676.8334
$ws.Range("I98").Value = 1874.7
$ws.Range("J98").Value = 687.5
$ws.Range("K98").Value = 1874.7
$ws.Range("L98").Value = 687.5
$ws.Range("M98").Value = -376.7
$ws.Range("N98").Value = -3683.5
$ws.Range("H122").Value = 1676.8334
$ws.Range("I122").Value = 1874.7
$ws.Range("J122").Value = 687.5
$ws.Range("K122").Value = 5624.1
$ws.Range("L122").Value = 2062.5
$ws.Range("M122").Value = -3174.1
$ws.Range("N122").Value = -6962.5
$ws.Range("H129").Value = 1719.7059
$ws.Range("I129").Value = 769.5833
$ws.Range("K129").Value = 2308.7499
$ws.Range("M129").Value = 2691.2501
$ws.Range("H137").Value = 1345004.6
$ws.Range("I137").Value = 65384.31
$ws.Range("J137").Value = 1918627.5
$ws.Range("K137").Value = 196152.93
$ws.Range("L137").Value = 5755882.5
$ws.Range("M137").Value = -193602.93
$ws.Range("N137").Value = -5760982.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15624878
$ws.Range("I32").Value = 16226652
$ws.Range("J32").Value = 10208908
$ws.Range("K32").Value = 16226652
$ws.Range("L32").Value = 10208908
$ws.Range("M32").Value = -16226365
$ws.Range("N32").Value = -10209482
$ws.Range("H37").Value = 52813.816
$ws.Range("J37").Value = 68554.86
$ws.Range("L37").Value = 68554.86
$ws.Range("N37").Value = -69100.86
$ws.Range("H45").Value = 4991.2607
$ws.Range("I45").Value = 5216.5835
$ws.Range("J45").Value = 4745.4546
$ws.Range("K45").Value = 5216.5835
$ws.Range("L45").Value = 4745.4546
$ws.Range("M45").Value = -4839.5835
$ws.Range("N45").Value = -5499.4546
$ws.Range("H69").Value = 99995
$ws.Range("J69").Value = 99995
$ws.Range("L69").Value = 99995
$ws.Range("N69").Value = -101493
$ws.Range("H72").Value = 99995
$ws.Range("J72").Value = 99995
$ws.Range("L72").Value = 299985
$ws.Range("N72").Value = -307473
$ws.Range("H74").Value = 3139.9333
$ws.Range("I74").Value = 3007.0715
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 3007.0715
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -2133.0715
$ws.Range("N74").Value = -6748
$ws.Range("H77").Value = 3139.9333
$ws.Range("I77").Value = 3007.0715
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 15035.3575
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -10667.3575
$ws.Range("N77").Value = -33736
$ws.Range("H132").Value = 2209.6086
$ws.Range("I132").Value = 2105.3777
$ws.Range("K132").Value = 6316.1331
$ws.Range("M132").Value = -3786.1331

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2780242.8
$ws.Range("I134").Value = 3511413.8
$ws.Range("J134").Value = 1793.2
$ws.Range("K134").Value = 10534241.4
$ws.Range("L134").Value = 5379.6
$ws.Range("M134").Value = -10531706.4
$ws.Range("N134").Value = -10449.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 642.6667
$ws.Range("I16").Value = 599.1177
$ws.Range("J16").Value = 716.7
$ws.Range("K16").Value = 599.1177
$ws.Range("L16").Value = 716.7
$ws.Range("M16").Value = -312.1177
$ws.Range("N16").Value = -1290.7
$ws.Range("H31").Value = 4534.915
$ws.Range("I31").Value = 2499
$ws.Range("J31").Value = 4777.2856
$ws.Range("K31").Value = 2499
$ws.Range("L31").Value = 4777.2856
$ws.Range("M31").Value = -2204
$ws.Range("N31").Value = -5367.2856
$ws.Range("H34").Value = 4534.915
$ws.Range("I34").Value = 2499
$ws.Range("J34").Value = 4777.2856
$ws.Range("K34").Value = 2499
$ws.Range("L34").Value = 4777.2856
$ws.Range("M34").Value = -2297
$ws.Range("N34").Value = -5181.2856
$ws.Range("H113").Value = 642.6667
$ws.Range("I113").Value = 599.1177
$ws.Range("J113").Value = 716.7
$ws.Range("K113").Value = 599.1177
$ws.Range("L113").Value = 716.7
$ws.Range("M113").Value = 1570.8823
$ws.Range("N113").Value = -5056.7
$ws.Range("H133").Value = 29888
$ws.Range("J133").Value = 29888
$ws.Range("L133").Value = 29888
$ws.Range("N133").Value = -34948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2574.55
$ws.Range("I68").Value = 2142.5715
$ws.Range("J68").Value = 2807.1538
$ws.Range("K68").Value = 6427.7145
$ws.Range("L68").Value = 8421.4614
$ws.Range("M68").Value = -5616.7145
$ws.Range("N68").Value = -10043.4614
$ws.Range("H71").Value = 2574.55
$ws.Range("I71").Value = 2142.5715
$ws.Range("J71").Value = 2807.1538
$ws.Range("K71").Value = 19283.1435
$ws.Range("L71").Value = 25264.3842
$ws.Range("M71").Value = -15227.1435
$ws.Range("N71").Value = -33376.3842
$ws.Range("H132").Value = 1114833.8
$ws.Range("I132").Value = 3666.6667
$ws.Range("J132").Value = 1670417.4
$ws.Range("K132").Value = 33000.0003
$ws.Range("L132").Value = 15033756.6
$ws.Range("M132").Value = -30470.0003
$ws.Range("N132").Value = -15038816.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 257500750
$ws.Range("I14").Value = 505000000
$ws.Range("J14").Value = 10001500
$ws.Range("K14").Value = 505000000
$ws.Range("L14").Value = 10001500
$ws.Range("M14").Value = -504999832
$ws.Range("N14").Value = -10001836
$ws.Range("H80").Value = 3600.5715
$ws.Range("I80").Value = 3740.8
$ws.Range("J80").Value = 3250
$ws.Range("K80").Value = 3740.8
$ws.Range("L80").Value = 3250
$ws.Range("M80").Value = -2742.8
$ws.Range("N80").Value = -5246
$ws.Range("H83").Value = 3600.5715
$ws.Range("I83").Value = 3740.8
$ws.Range("J83").Value = 3250
$ws.Range("K83").Value = 18704
$ws.Range("L83").Value = 16250
$ws.Range("M83").Value = -13712
$ws.Range("N83").Value = -26234
$ws.Range("H102").Value = 4995.8335
$ws.Range("I102").Value = 5494
$ws.Range("K102").Value = 5494
$ws.Range("M102").Value = -3872
$ws.Range("H132").Value = 2265.2
$ws.Range("I132").Value = 2175.5
$ws.Range("J132").Value = 2624
$ws.Range("K132").Value = 6526.5
$ws.Range("L132").Value = 7872
$ws.Range("M132").Value = -3996.5
$ws.Range("N132").Value = -12932

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3294
$ws.Range("I22").Value = 2243.4285
$ws.Range("J22").Value = 4213.25
$ws.Range("K22").Value = 2243.4285
$ws.Range("L22").Value = 4213.25
$ws.Range("M22").Value = -1948.4285
$ws.Range("N22").Value = -4803.25
$ws.Range("H27").Value = 3294
$ws.Range("I27").Value = 2243.4285
$ws.Range("J27").Value = 4213.25
$ws.Range("K27").Value = 2243.4285
$ws.Range("L27").Value = 4213.25
$ws.Range("M27").Value = -2136.4285
$ws.Range("N27").Value = -4427.25
$ws.Range("H100").Value = 4725.75
$ws.Range("I100").Value = 4667.6665
$ws.Range("J100").Value = 4900
$ws.Range("K100").Value = 4667.6665
$ws.Range("L100").Value = 4900
$ws.Range("M100").Value = -4126.6665
$ws.Range("N100").Value = -5982
$ws.Range("H122").Value = 14224.917
$ws.Range("I122").Value = 17783.834
$ws.Range("J122").Value = 10666
$ws.Range("K122").Value = 53351.50199999999
$ws.Range("L122").Value = 31998
$ws.Range("M122").Value = -50901.50199999999
$ws.Range("N122").Value = -36898
$ws.Range("H139").Value = 67714.664
$ws.Range("J139").Value = 67714.664
$ws.Range("L139").Value = 67714.664
$ws.Range("N139").Value = -77994.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2624.6667
$ws.Range("I132").Value = 2338.7222
$ws.Range("J132").Value = 3482.5
$ws.Range("K132").Value = 7016.1666
$ws.Range("L132").Value = 10447.5
$ws.Range("M132").Value = -4486.1666
$ws.Range("N132").Value = -15507.5
